$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab (was "Sheet1")
$ws.Name = "Surrogate Detection Statistics"

# Unhide every column on the sheet - all of the analysis "helper" columns
# that used to be hidden are shown again.
$ws.Cells.EntireColumn.Hidden = $false

# Re-apply each column's (slightly relaid-out) width. Columns that keep an
# auto "best fit" width on both sides of the edit (D and the BY:CC summary
# columns) are intentionally left alone so Excel keeps auto-sizing them.
$columnWidths = [ordered]@{
    "A:A"   = 8.690104166666666
    "B:B"   = 9.690104166666666
    "C:C"   = 10.955729166666666
    "E:E"   = 14.061197916666666
    "F:F"   = 13.795572916666666
    "G:G"   = 7.955729166666667
    "H:H"   = 12.901041666666666
    "I:I"   = 14.795572916666666
    "J:J"   = 12.901041666666666
    "K:K"   = 21.533854166666668
    "L:L"   = 22.166666666666668
    "M:M"   = 11.061197916666666
    "N:O"   = 7.955729166666667
    "P:P"   = 8.268229166666666
    "Q:Q"   = 15.795572916666666
    "R:R"   = 15.061197916666666
    "S:S"   = 16.901041666666668
    "T:T"   = 30.639322916666668
    "U:U"   = 18.639322916666668
    "V:V"   = 22.955729166666668
    "W:W"   = 19.795572916666668
    "X:AD"  = 8.955729166666666
    "AE:AI" = 9.955729166666666
    "AJ:AJ" = 8.955729166666666
    "AK:AK" = 7.955729166666667
    "AL:AL" = 8.955729166666666
    "AM:AM" = 7.428385416666667
    "AN:AN" = 8.955729166666666
    "AO:AO" = 11.533854166666666
    "AP:AR" = 12.639322916666666
    "AS:AS" = 13.639322916666666
    "AT:AT" = 8.061197916666666
    "AU:AU" = 9.901041666666666
    "AV:AX" = 10.901041666666666
    "AY:AY" = 11.901041666666666
    "AZ:AZ" = 8.955729166666666
    "BA:BC" = 9.955729166666666
    "BD:BD" = 10.955729166666666
    "BE:BE" = 20.690104166666668
    "BF:BH" = 21.690104166666668
    "BI:BI" = 22.690104166666668
    "BJ:BJ" = 25.061197916666668
    "BK:BM" = 26.061197916666668
    "BN:BN" = 27.061197916666668
    "BO:BO" = 19.901041666666668
    "BP:BQ" = 20.955729166666668
    "BR:BR" = 21.795572916666668
    "BS:BS" = 21.955729166666668
    "BT:BT" = 10.795572916666666
    "BU:BW" = 11.795572916666666
    "BX:BX" = 12.901041666666666
}
foreach ($range in $columnWidths.Keys) {
    $ws.Columns($range).ColumnWidth = $columnWidths[$range]
}

# Scroll the frozen pane back to the top (topLeftCell A2) and move the
# active selection in the lower pane to H20.
$ws.Range("A2").Select() | Out-Null
$ws.Range("H20").Select() | Out-Null
